$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values are written as literal text (matching the source inlineStr cells)
# rather than being auto-parsed into numbers/percentages by Excel.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "329.22"
$ws.Range("E2").Value = "-0.64%"
$ws.Range("D3").Value = "43.54"
$ws.Range("E3").Value = "4.10%"
$ws.Range("D4").Value = "5.562"
$ws.Range("E4").Value = "-2.03%"
$ws.Range("D5").Value = "0.08207"
$ws.Range("E5").Value = "-2.11%"
$ws.Range("D6").Value = "8.760"
$ws.Range("E6").Value = "-0.43%"
$ws.Range("D7").Value = "4.354"
$ws.Range("E7").Value = "-3.87%"
$ws.Range("D8").Value = "1.879"
$ws.Range("E8").Value = "-6.84%"
$ws.Range("D9").Value = "2.838"
$ws.Range("E9").Value = "-4.60%"
$ws.Range("D10").Value = "0.9440"
$ws.Range("E10").Value = "1.77%"
$ws.Range("D11").Value = "0.1190"
$ws.Range("E11").Value = "-8.40%"
$ws.Range("D12").Value = "0.1907"
$ws.Range("E12").Value = "-3.46%"
$ws.Range("D13").Value = "0.09703"
$ws.Range("E13").Value = "2.74%"
$ws.Range("E14").Value = "10.94%"
$ws.Range("D15").Value = "0.1069"
$ws.Range("E15").Value = "0.82%"
$ws.Range("D16").Value = "0.001284"
$ws.Range("E16").Value = "-1.13%"
$ws.Range("D17").Value = "0.005960"
$ws.Range("E17").Value = "-2.06%"
$ws.Range("D18").Value = "3.532"
$ws.Range("E18").Value = "2.82%"
$ws.Range("D20").Value = "8.744"
$ws.Range("E20").Value = "8.77%"
$ws.Range("E21").Value = "-0.19%"
$ws.Range("D22").Value = "0.2499"
$ws.Range("E22").Value = "-4.31%"
$ws.Range("D23").Value = "0.04404"
$ws.Range("E23").Value = "-0.74%"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").Value = "-1.42%"
$ws.Range("D25").Value = "0.004303"
$ws.Range("E25").Value = "-1.68%"
$ws.Range("E26").Value = "2.89%"
$ws.Range("E27").Value = "31.68%"
$ws.Range("D39").Value = "0.02741"
$ws.Range("E39").Value = "-2.17%"
$ws.Range("D40").Value = "0.05687"
$ws.Range("E40").Value = "2.83%"
$ws.Range("D41").Value = "0.007874"
$ws.Range("E41").Value = "0.85%"
$ws.Range("D42").Value = "0.009747"
$ws.Range("E43").Value = "-1.14%"
$ws.Range("D44").Value = "0.002107"
$ws.Range("E44").Value = "-2.48%"
$ws.Range("D45").Value = "0.01002"
$ws.Range("E45").Value = "-8.14%"
$ws.Range("E46").Value = "4.40%"
$ws.Range("E47").Value = "0.37%"
$ws.Range("D48").Value = "0.003452"
$ws.Range("E48").Value = "-2.09%"
$ws.Range("E49").Value = "0.05%"
$ws.Range("E50").Value = "0.37%"
$ws.Range("E51").Value = "0.37%"

# Drop the temporary Text number format so the cells fall back to the
# workbook default style, matching the unstyled cells in the source file.
$textRange.ClearFormats()
